# Daily attendance processing - refreshes recorder lists and attendance
# counts/percentages on the "Session Analysis Results" sheet after the
# latest sync from the attendance system.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PercentText($CellAddress, $PercentText) {
    # Keep these as literal text (e.g. "60.2%") instead of letting Excel
    # reinterpret the assignment as a numeric percentage value.
    $r = $ws.Range($CellAddress)
    $r.NumberFormat = "@"
    $r.Value = $PercentText
}

# --- Year 2 / A1 - ANATOMY session 1 ---------------------------------------
$ws.Range("G2").Value = "nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("H2").Value = "187/216"

# --- Year 2 / A1 - ANATOMY session 2 ---------------------------------------
$ws.Range("G3").Value = "mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H3").Value = "143/216"

# --- Recomputed summary / class statistics percentages ---------------------
Set-PercentText "L10" "60.2%"
Set-PercentText "S15" "80.1%"
Set-PercentText "S16" "61.8%"

# --- Year 2 / A2 - ANATOMY session 1 ---------------------------------------
$ws.Range("G17").Value = "nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("H17").Value = "163/217"

# --- Year 2 / A2 - ANATOMY session 2 ---------------------------------------
$ws.Range("G18").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H18").Value = "122/217"

Set-PercentText "S18" "73.1%"

# --- Year 2 / A3 - ANATOMY session 1 ---------------------------------------
$ws.Range("G32").Value = "gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Year 2 / A3 - ANATOMY session 2 ---------------------------------------
$ws.Range("G33").Value = "gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"

# --- Year 2 / A4 - ANATOMY session 1 ---------------------------------------
$ws.Range("G47").Value = "gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"
$ws.Range("H47").Value = "150/225"

# --- Year 2 / A4 - ANATOMY session 2 ---------------------------------------
$ws.Range("G48").Value = "gehanadel@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, System, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg"
$ws.Range("H48").Value = "179/225"

# --- Year 2 / B1 - ANATOMY session 1 ---------------------------------------
$ws.Range("G62").Value = "gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Year 2 / B1 - ANATOMY session 2 ---------------------------------------
$ws.Range("G63").Value = "mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Year 2 / B2 - ANATOMY session 1 ---------------------------------------
$ws.Range("G77").Value = "gehanadel@med.asu.edu.eg, nahla.nagiub@med.asu.edu.eg, mennatulla.medhat@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, hend_mahmoud@med.asu.edu.eg, nourhan.mahmoud@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Year 2 / B2 - ANATOMY session 2 ---------------------------------------
$ws.Range("G78").Value = "asmaa.reda@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg"

# --- Year 2 / B3 - ANATOMY session 1 ---------------------------------------
$ws.Range("G92").Value = "nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"

# --- Year 2 / B4 - ANATOMY session 1 ---------------------------------------
$ws.Range("G107").Value = "nahla.nagiub@med.asu.edu.eg, rana.abozaid@med.asu.edu.eg, nesmadrahim@med.asu.edu.eg, servinaz@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
